$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-4: vervoerregio column (A) becomes an empty/zero numeric
# placeholder with default formatting; gewest column (B) keeps its
# numeric style but two of the reference values change. ---
$ws.Range("A2").Value = 0
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 2000

$ws.Range("A3").Value = 0
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 4000

$ws.Range("A4").Value = 0
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = 99999

# --- Rows 5-19: vervoerregio ids re-entered as text, sorted the way
# text sorts ("1","10","11","12","13","14","15","2","3",...,"9") — this
# is the new gemeente_gewest aggregate table, now carrying the 3 extra
# reference regions (7,8,9) appended as rows 17-19. ---
$values = @("1","10","11","12","13","14","15","2","3","4","5","6","7","8","9")

# Make sure the new rows exist with the same numeric style as the
# existing gewest column before we start writing into them.
$ws.Range("B16").Copy()
$ws.Range("B17:B19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 1).Value = "'" + $values[$i]
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Value = 2000
}
